$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the date in A1 by one day
$ws.Range("A1").Value = 45309

# Step 2: update price column (D14:D21) with recalculated values
$ws.Range("D14").Value = 76.307
$ws.Range("D15").Value = 99.298
$ws.Range("D16").Value = 118.919
$ws.Range("D17").Value = 189.676
$ws.Range("D18").Value = 236.649
$ws.Range("D19").Value = 295.514
$ws.Range("D20").Value = 330.991
$ws.Range("D21").Value = 366.666
